{"js": "// This edit renames the header date and swaps every two-digit-by-two-digit\n// multiplication prompt in the practice-sheet table for a new value.\n// Each \"old\" value below occurs exactly once in the document, so a plain\n// search-and-replace per pair is unambiguous and leaves every run's\n// original formatting (font/size) untouched.\nconst replacements = [\n  [\"2025-10-04 Saturday\", \"2025-10-05 Sunday\"],\n  [\"40\u00d731=\", \"88\u00d762=\"],\n  [\"75\u00d749=\", \"81\u00d753=\"],\n  [\"22\u00d784=\", \"90\u00d791=\"],\n  [\"48\u00d719=\", \"49\u00d748=\"],\n  [\"82\u00d775=\", \"43\u00d762=\"],\n  [\"70\u00d775=\", \"15\u00d796=\"],\n  [\"37\u00d798=\", \"62\u00d760=\"],\n  [\"34\u00d791=\", \"39\u00d752=\"],\n  [\"67\u00d715=\", \"45\u00d769=\"],\n  [\"59\u00d762=\", \"44\u00d785=\"],\n  [\"94\u00d779=\", \"93\u00d788=\"],\n  [\"80\u00d748=\", \"20\u00d776=\"],\n  [\"66\u00d723=\", \"54\u00d772=\"],\n  [\"17\u00d752=\", \"74\u00d762=\"],\n  [\"94\u00d714=\", \"92\u00d788=\"],\n  [\"89\u00d741=\", \"53\u00d743=\"],\n  [\"44\u00d766=\", \"58\u00d782=\"],\n  [\"95\u00d765=\", \"19\u00d762=\"],\n  [\"37\u00d768=\", \"47\u00d747=\"],\n  [\"88\u00d782=\", \"79\u00d775=\"],\n  [\"49\u00d747=\", \"16\u00d731=\"],\n  [\"40\u00d766=\", \"32\u00d715=\"],\n  [\"30\u00d736=\", \"82\u00d798=\"],\n  [\"14\u00d781=\", \"85\u00d720=\"],\n  [\"86\u00d732=\", \"27\u00d737=\"]\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and every two-digit-by-two-digit multiplication\n# expression cell with its new value. Every \"Old\" string below is unique\n# in the document, so a straightforward Find/Replace per pair is safe and\n# keeps each run's original formatting (font/size) untouched.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"2025-10-04 Saturday\"; New = \"2025-10-05 Sunday\" },\n    @{ Old = \"40\u00d731=\"; New = \"88\u00d762=\" },\n    @{ Old = \"75\u00d749=\"; New = \"81\u00d753=\" },\n    @{ Old = \"22\u00d784=\"; New = \"90\u00d791=\" },\n    @{ Old = \"48\u00d719=\"; New = \"49\u00d748=\" },\n    @{ Old = \"82\u00d775=\"; New = \"43\u00d762=\" },\n    @{ Old = \"70\u00d775=\"; New = \"15\u00d796=\" },\n    @{ Old = \"37\u00d798=\"; New = \"62\u00d760=\" },\n    @{ Old = \"34\u00d791=\"; New = \"39\u00d752=\" },\n    @{ Old = \"67\u00d715=\"; New = \"45\u00d769=\" },\n    @{ Old = \"59\u00d762=\"; New = \"44\u00d785=\" },\n    @{ Old = \"94\u00d779=\"; New = \"93\u00d788=\" },\n    @{ Old = \"80\u00d748=\"; New = \"20\u00d776=\" },\n    @{ Old = \"66\u00d723=\"; New = \"54\u00d772=\" },\n    @{ Old = \"17\u00d752=\"; New = \"74\u00d762=\" },\n    @{ Old = \"94\u00d714=\"; New = \"92\u00d788=\" },\n    @{ Old = \"89\u00d741=\"; New = \"53\u00d743=\" },\n    @{ Old = \"44\u00d766=\"; New = \"58\u00d782=\" },\n    @{ Old = \"95\u00d765=\"; New = \"19\u00d762=\" },\n    @{ Old = \"37\u00d768=\"; New = \"47\u00d747=\" },\n    @{ Old = \"88\u00d782=\"; New = \"79\u00d775=\" },\n    @{ Old = \"49\u00d747=\"; New = \"16\u00d731=\" },\n    @{ Old = \"40\u00d766=\"; New = \"32\u00d715=\" },\n    @{ Old = \"30\u00d736=\"; New = \"82\u00d798=\" },\n    @{ Old = \"14\u00d781=\"; New = \"85\u00d720=\" },\n    @{ Old = \"86\u00d732=\"; New = \"27\u00d737=\" }\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair.Old\n    $find.Replacement.Text = $pair.New\n    $find.Execute($pair.Old, $false, $false, $false, $false, $false, $true, 0, $false, $pair.New, 2) | Out-Null\n}\n"}
